# Applies the cryptos-list refresh described in the commit diff:
# updates the "Price" (D) and "Volume(1h)" (E) columns for rows 2-51.
#
# Many "Price" values look like plain numbers/dates to Excel's auto-detection
# (e.g. "1.009", "216.41"), which would silently coerce them into numeric
# cells if written with a plain .Value assignment. The source file stores
# them as text, so for every such value we stage it in a scratch cell that
# has been forced to Text format, copy it, and paste-special (values only)
# into the destination. That keeps the destination cell text-typed without
# permanently tattooing a NumberFormat style onto it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = "26.218.59"
$ws.Cells.Item(2, 5).Value = "  -4.15%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.655.57"
$ws.Cells.Item(3, 5).Value = "  -3.53%  "

# Row 4
$scratch.Value = "1.008"
$scratch.Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4163)
$ws.Cells.Item(4, 5).Value = "  +0.13%  "

# Row 5
$scratch.Value = "216.41"
$scratch.Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163)
$ws.Cells.Item(5, 5).Value = "  -3.79%  "

# Row 6
$scratch.Value = "0.5138"
$scratch.Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163)
$ws.Cells.Item(6, 5).Value = "  -3.04%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.20%  "

# Row 8
$scratch.Value = "0.2596"
$scratch.Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4163)
$ws.Cells.Item(8, 5).Value = "  -2.15%  "

# Row 9
$scratch.Value = "0.06450"
$scratch.Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4163)
$ws.Cells.Item(9, 5).Value = "  -3.55%  "

# Row 10
$scratch.Value = "19.85"
$scratch.Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4163)
$ws.Cells.Item(10, 5).Value = "  -5.22%  "

# Row 11
$scratch.Value = "0.07835"
$scratch.Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4163)

# Row 12
$ws.Cells.Item(12, 4).Value = "1.657.88"
$ws.Cells.Item(12, 5).Value = "  -3.08%  "

# Row 13
$scratch.Value = "4.301"
$scratch.Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4163)
$ws.Cells.Item(13, 5).Value = "  -4.20%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "1.884.35"
$ws.Cells.Item(14, 5).Value = "  -3.51%  "

# Row 15
$scratch.Value = "0.5524"
$scratch.Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4163)
$ws.Cells.Item(15, 5).Value = "  -4.69%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "0.0₅8030"
$ws.Cells.Item(16, 5).Value = "  -1.95%  "

# Row 17
$scratch.Value = "64.19"
$scratch.Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4163)
$ws.Cells.Item(17, 5).Value = "  -5.28%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "26.243.96"
$ws.Cells.Item(18, 5).Value = "  -4.14%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  +0.18%  "

# Row 20
$scratch.Value = "210.53"
$scratch.Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4163)
$ws.Cells.Item(20, 5).Value = "  -4.52%  "

# Row 21
$scratch.Value = "4.412"
$scratch.Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4163)
$ws.Cells.Item(21, 5).Value = "  -5.11%  "

# Row 22
$scratch.Value = "10.09"
$scratch.Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163)
$ws.Cells.Item(22, 5).Value = "  -3.34%  "

# Row 23
$scratch.Value = "6.024"
$scratch.Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4163)
$ws.Cells.Item(23, 5).Value = "  -0.11%  "

# Row 25
$scratch.Value = "144.93"
$scratch.Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4163)
$ws.Cells.Item(25, 5).Value = "  -0.53%  "

# Row 26
$scratch.Value = "1.795"
$scratch.Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4163)
$ws.Cells.Item(26, 5).Value = "  +5.05%  "

# Row 27
$scratch.Value = "0.1175"
$scratch.Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163)
$ws.Cells.Item(27, 5).Value = "  -2.69%  "

# Row 28
$scratch.Value = "7.009"
$scratch.Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4163)
$ws.Cells.Item(28, 5).Value = "  -3.31%  "

# Row 29
$scratch.Value = "15.86"
$scratch.Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163)
$ws.Cells.Item(29, 5).Value = "  -2.03%  "

# Row 30
$scratch.Value = "0.05110"
$scratch.Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4163)
$ws.Cells.Item(30, 5).Value = "  -5.01%  "

# Row 31
$scratch.Value = "1.243"
$scratch.Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163)
$ws.Cells.Item(31, 5).Value = "  -4.15%  "

# Row 32
$scratch.Value = "3.364"
$scratch.Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4163)
$ws.Cells.Item(32, 5).Value = "  -3.43%  "

# Row 33
$scratch.Value = "3.237"
$scratch.Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4163)
$ws.Cells.Item(33, 5).Value = "  -4.75%  "

# Row 34
$scratch.Value = "1.561"
$scratch.Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4163)
$ws.Cells.Item(34, 5).Value = "  -4.66%  "

# Row 35
$scratch.Value = "2.737"
$scratch.Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4163)
$ws.Cells.Item(35, 5).Value = "  -4.12%  "

# Row 36
$scratch.Value = "2.358"
$scratch.Copy() | Out-Null
$ws.Range("D36").PasteSpecial(-4163)
$ws.Cells.Item(36, 5).Value = "  -1.63%  "

# Row 37
$scratch.Value = "0.9210"
$scratch.Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4163)

# Row 38
$scratch.Value = "0.5726"
$scratch.Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163)
$ws.Cells.Item(38, 5).Value = "  -2.77%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "1.167.26"
$ws.Cells.Item(39, 5).Value = "  +0.75%  "

# Row 40
$scratch.Value = "0.01590"
$scratch.Copy() | Out-Null
$ws.Range("D40").PasteSpecial(-4163)
$ws.Cells.Item(40, 5).Value = "  -3.74%  "

# Row 41
$scratch.Value = "2.563"
$scratch.Copy() | Out-Null
$ws.Range("D41").PasteSpecial(-4163)
$ws.Cells.Item(41, 5).Value = "  -0.60%  "

# Row 42
$scratch.Value = "1.009"
$scratch.Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4163)
$ws.Cells.Item(42, 5).Value = "  +0.17%  "

# Row 43
$scratch.Value = "5.694"
$scratch.Copy() | Out-Null
$ws.Range("D43").PasteSpecial(-4163)
$ws.Cells.Item(43, 5).Value = "  -2.39%  "

# Row 44
$scratch.Value = "0.8289"
$scratch.Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4163)
$ws.Cells.Item(44, 5).Value = "  -1.38%  "

# Row 45
$scratch.Value = "100.34"
$scratch.Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4163)
$ws.Cells.Item(45, 5).Value = "  -0.68%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "1.797.95"
$ws.Cells.Item(46, 5).Value = "  -3.31%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "0.0₈112"
$ws.Cells.Item(47, 5).Value = "  -4.94%  "

# Row 48
$scratch.Value = "0.4546"
$scratch.Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4163)
$ws.Cells.Item(48, 5).Value = "  -0.52%  "

# Row 49
$scratch.Value = "55.45"
$scratch.Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4163)
$ws.Cells.Item(49, 5).Value = "  -4.09%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  +0.17%  "

# Row 51
$scratch.Value = "7.893"
$scratch.Copy() | Out-Null
$ws.Range("D51").PasteSpecial(-4163)
$ws.Cells.Item(51, 5).Value = "  -3.36%  "

$scratch.Clear()
$excel.CutCopyMode = 0
